# Updated cryptos list on Wed Apr 19 05:39:12 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price ("D") cells whose new text would otherwise be auto-recognized by Excel
# as a genuine number (losing their original plain-text storage / exact
# formatting). Force them to remain text via the "@" (Text) number format
# before writing the value, exactly like the source cells already are.
$textPriceCells = @(
    "D5","D7","D8","D9","D12","D15","D16","D17","D18","D19","D20","D21","D22",
    "D26","D27","D28","D29","D30","D32","D34","D35","D36","D37","D38","D39",
    "D40","D41","D43","D44","D45","D47","D48","D50"
)
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "30.283.14"
$ws.Range("E2").Value = "  +1.71%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.091.01"
$ws.Range("E3").Value = "  -0.63%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.48%  "

# Row 5 - BNB
$ws.Range("D5").Value = "341.97"
$ws.Range("E5").Value = "  -1.62%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  -0.49%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.5310"
$ws.Range("E7").Value = "  +2.11%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.4387"
$ws.Range("E8").Value = "  -0.78%  "

# Row 9 - OKB
$ws.Range("D9").Value = "54.33"
$ws.Range("E9").Value = "  +0.66%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -0.39%  "

# Row 11 - Polygon
$ws.Range("E11").Value = "  +0.02%  "

# Row 12 - Solana
$ws.Range("D12").Value = "24.69"
$ws.Range("E12").Value = "  -0.93%  "

# Row 13 - Chainlink
$ws.Range("E13").Value = "  +3.32%  "

# Row 14 / 15 - Polkadot and WrappedEther swapped places
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.981.90"
$ws.Range("E14").Value = "  -4.14%  "

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "6.885"
$ws.Range("E15").Value = "  +0.95%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "101.65"
$ws.Range("E16").Value = "  -1.12%  "

# Row 17 - ShibaInu
$ws.Range("D17").Value = "0.00001160"
$ws.Range("E17").Value = "  -0.02%  "

# Row 18 - BinanceUSD
$ws.Range("D18").Value = "1.002"
$ws.Range("E18").Value = "  -0.64%  "

# Row 19 - Avalanche
$ws.Range("D19").Value = "21.15"
$ws.Range("E19").Value = "  -0.07%  "

# Row 20 - TRON
$ws.Range("D20").Value = "0.06716"
$ws.Range("E20").Value = "  +0.62%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "6.345"
$ws.Range("E21").Value = "  +0.74%  "

# Row 22 - Dai
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  -0.48%  "

# Row 23 - WrappedBTC
$ws.Range("D23").Value = "30.306.56"
$ws.Range("E23").Value = "  +1.72%  "

# Row 24 - Cosmos
$ws.Range("E24").Value = "  -1.25%  "

# Row 25 - Toncoin
$ws.Range("E25").Value = "  -0.43%  "

# Row 26 - EthereumClassic
$ws.Range("D26").Value = "21.81"
$ws.Range("E26").Value = "  -0.97%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("D27").Value = "6.883"
$ws.Range("E27").Value = "  +8.25%  "

# Row 28 - Monero
$ws.Range("D28").Value = "162.88"
$ws.Range("E28").Value = "  +0.29%  "

# Row 29 - LidoDAOToken
$ws.Range("D29").Value = "2.497"
$ws.Range("E29").Value = "  -1.21%  "

# Row 30 - BitcoinCash
$ws.Range("D30").Value = "133.65"
$ws.Range("E30").Value = "  +0.02%  "

# Row 31 - ImmutableX
$ws.Range("E31").Value = "  -0.73%  "

# Row 32 - ARBITRUM
$ws.Range("D32").Value = "1.669"
$ws.Range("E32").Value = "  -6.38%  "

# Row 33 - Stellar
$ws.Range("E33").Value = "  -0.31%  "

# Row 34 - Filecoin
$ws.Range("D34").Value = "6.282"
$ws.Range("E34").Value = "  +0.93%  "

# Row 35 - HuobiToken
$ws.Range("D35").Value = "3.912"
$ws.Range("E35").Value = "  -0.88%  "

# Row 36 - FraxShare
$ws.Range("D36").Value = "10.13"
$ws.Range("E36").Value = "  -3.53%  "

# Row 37 - VeChain
$ws.Range("D37").Value = "0.02618"
$ws.Range("E37").Value = "  +1.12%  "

# Row 38 - Hedera
$ws.Range("D38").Value = "0.06763"
$ws.Range("E38").Value = "  -0.07%  "

# Row 39 - Aptos
$ws.Range("D39").Value = "12.61"
$ws.Range("E39").Value = "  +0.03%  "

# Row 40 - TheSandbox
$ws.Range("D40").Value = "0.6963"
$ws.Range("E40").Value = "  -0.60%  "

# Row 41 - TrustWalletToken
$ws.Range("D41").Value = "1.342"
$ws.Range("E41").Value = "  +0.79%  "

# Row 42 - Algorand
$ws.Range("E42").Value = "  -0.81%  "

# Row 43 - Decentraland
$ws.Range("D43").Value = "0.6756"
$ws.Range("E43").Value = "  -1.12%  "

# Row 44 - NEARProtocol
$ws.Range("D44").Value = "2.383"
$ws.Range("E44").Value = "  +1.15%  "

# Row 45 - EnergySwap
$ws.Range("D45").Value = "14.29"
$ws.Range("E45").Value = "  -0.91%  "

# Row 46 - Frax
$ws.Range("E46").Value = "  -0.43%  "

# Row 47 - WEMIXTOKEN
$ws.Range("D47").Value = "1.282"
$ws.Range("E47").Value = "  +5.57%  "

# Row 48 - PancakeSwap
$ws.Range("D48").Value = "3.632"
$ws.Range("E48").Value = "  -0.24%  "

# Row 49 - BabyDogeCoin
$ws.Range("E49").Value = "  -2.31%  "

# Row 50 - ThetaToken
$ws.Range("D50").Value = "1.206"
$ws.Range("E50").Value = "  +2.82%  "

# Row 51 - EOS
$ws.Range("E51").Value = "  -0.79%  "
